$d = $word.ActiveDocument

$pairs = @(
    @("775×6=", "921×6="),
    @("483×2=", "358×4="),
    @("277×7=", "486×3="),
    @("288×8=", "423×2="),
    @("392×7=", "883×6="),
    @("580×2=", "759×9="),
    @("412×5=", "129×9="),
    @("725×3=", "767×3="),
    @("531×5=", "233×2="),
    @("425×4=", "404×8="),
    @("253×5=", "211×9="),
    @("638×2=", "528×7="),
    @("563×3=", "927×3="),
    @("608×6=", "749×3="),
    @("519×6=", "522×6="),
    @("523×4=", "215×7="),
    @("136×3=", "275×7="),
    @("372×2=", "435×4="),
    @("685×7=", "731×6="),
    @("489×3=", "142×9="),
    @("547×7=", "919×3="),
    @("115×7=", "116×8="),
    @("596×2=", "381×7="),
    @("445×9=", "783×9="),
    @("211×4=", "240×5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
